$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Template cells to copy number-format/style from (row 42 already has the
# date style in column B and the wrapped-text style in column D).
$dateTemplate = $ws.Cells.Item(42, 2)
$textTemplate = $ws.Cells.Item(42, 4)

function Add-LogRow {
    param($Row, $DateSerial, $Task, $Status, $TallRow)

    $dateTemplate.Copy()
    $ws.Cells.Item($Row, 2).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($Row, 2).Value = $DateSerial

    if ($Task) {
        $textTemplate.Copy()
        $ws.Cells.Item($Row, 4).PasteSpecial(-4122) | Out-Null
        $ws.Cells.Item($Row, 4).Value = $Task
    }

    if ($Status) {
        $ws.Cells.Item($Row, 5).Value = $Status
    }

    if ($TallRow) {
        $ws.Rows.Item($Row).RowHeight = 30
    }
}

Add-LogRow 43 43167 "Created forms and tried to push values to database through forms" "completed" $true
Add-LogRow 44 43168 "automated data update in database through forms" "completed" $true
Add-LogRow 45 43170 "Found bugs in data entry and trying to resolve it" "in progress" $true
Add-LogRow 46 43171 "Found bugs in data entry and trying to resolve it" "in progress" $true
Add-LogRow 47 43172 "Found bugs in data entry and trying to resolve it" "in progress" $true
Add-LogRow 48 43173 "Found alternative method for pushing data to database through build data table" "completed" $true
Add-LogRow 49 43174 "Implemented build datatable method for data pushing" "completed" $true
Add-LogRow 50 43175 "Found errors in pushing data" "in progress" $false
Add-LogRow 51 43176 $null $null $false
Add-LogRow 52 43177 $null $null $false
Add-LogRow 53 43178 "Trying to resolve bugs" "in progress" $false
Add-LogRow 54 43179 "Resolved bugs" "completed" $false
Add-LogRow 55 43180 "Integrating project and debugging" $null $false
Add-LogRow 56 43181 "Integrating project and debugging" $null $false
Add-LogRow 57 43182 "Integrating project and debugging" $null $false
Add-LogRow 58 43183 "Integrating project and debugging" $null $false
Add-LogRow 59 43184 "Integrating project and debugging" $null $false
Add-LogRow 60 43185 "Integrating project and debugging" $null $false

# Move the view/selection to reflect the newly added rows, mirroring the
# author's final cursor position after entering the data.
$ws.Activate() | Out-Null
$ws.Range("A60").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 46
$excel.ActiveWindow.ScrollColumn = 1
